# "Generate Report for Handoff"
#
# This localization-status report is regenerated: the status moves from
# "In Translation" to "Ready for handoff", the two "Latest ... Datetime"
# timestamps advance by about a minute (handoff just completed), and the
# "Status"/timestamp columns on each sheet widen to fit the new, longer
# status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps advance a little more than a minute -------------------------
# Overview!G2 and de-de!H2 shared the same "Latest HO Xliff Generate Date" /
# "Latest Handoff Datetime" timestamp.
$wsOverview.Range("G2").Value = "2016-10-20 09:27:43"
$wsDeDe.Range("H2").Value = "2016-10-20 09:27:43"

# zh-cn!H2 "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-10-20 09:27:31"

# --- Widen the Status / Datetime columns to fit "Ready for handoff" --------
# These columns grow from a stored width of ~13.41 to ~17.22 characters,
# which corresponds to a COM ColumnWidth of about 16.33.
$newColumnWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
